$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 56, shifting existing rows 56-64 down to 57-65.
$ws.Rows.Item(56).EntireRow.Insert()

# Populate the newly inserted row 56 with the new weekly record
# (same dimension values as the former last row, 7 days later).
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C56").Value = "Los Lagos"
$ws.Range("D56").Value = 44918
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = "Frutos de hueso (carozo)"
$ws.Range("I56").Value = 100103003
$ws.Range("J56").Value = "Damasco"
$ws.Range("K56").Value = "Castle Brite"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 600
$ws.Range("N56").Value = 19000
$ws.Range("O56").Value = 20000
$ws.Range("P56").Value = 19500
$ws.Range("Q56").Value = "$/caja 16 kilos"
$ws.Range("R56").Value = "Región de O'Higgins"
$ws.Range("S56").Value = 1219
$ws.Range("T56").Value = 16
